# Updated cryptos list with fresh price/volume data (GitHub Actions scrape).
# Column D/E are kept as text (NumberFormat '@') to match the source data's
# inline-string representation, e.g. '1.012' stays text instead of becoming 1.012.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '20.745.66'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -5.82%  '
$ws.Range('E2').Style = "Normal"

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.481.45'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -4.66%  '
$ws.Range('E3').Style = "Normal"

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.012'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +1.13%  '
$ws.Range('E4').Style = "Normal"

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.007'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('E5').Style = "Normal"

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '279.00'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -3.80%  '
$ws.Range('E6').Style = "Normal"

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3773'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -4.91%  '
$ws.Range('E7').Style = "Normal"

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3077'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -4.26%  '
$ws.Range('E8').Style = "Normal"

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '42.25'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -3.63%  '
$ws.Range('E9').Style = "Normal"

# Row 10
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.014'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -5.64%  '
$ws.Range('E10').Style = "Normal"

# Row 11
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06693'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -7.30%  '
$ws.Range('E11').Style = "Normal"

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.013'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('E12').Style = "Normal"

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.465'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('E13').Style = "Normal"

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.57'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -5.81%  '
$ws.Range('E14').Style = "Normal"

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.489.61'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -4.02%  '
$ws.Range('E15').Style = "Normal"

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.267'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -5.25%  '
$ws.Range('E16').Style = "Normal"

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001032'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -8.34%  '
$ws.Range('E17').Style = "Normal"

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06509'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('E18').Style = "Normal"

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '80.31'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -3.72%  '
$ws.Range('E19').Style = "Normal"

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.008'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('E20').Style = "Normal"

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.864'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -6.09%  '
$ws.Range('E21').Style = "Normal"

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.75'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -4.74%  '
$ws.Range('E22').Style = "Normal"

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.69'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -5.28%  '
$ws.Range('E23').Style = "Normal"

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.337'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('E24').Style = "Normal"

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '20.747.00'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -5.86%  '
$ws.Range('E25').Style = "Normal"

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.280'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -5.39%  '
$ws.Range('E26').Style = "Normal"

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '145.45'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('E27').Style = "Normal"

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.68'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -4.57%  '
$ws.Range('E28').Style = "Normal"

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.769'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.27%  '
$ws.Range('E29').Style = "Normal"

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.661.65'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -3.75%  '
$ws.Range('E30').Style = "Normal"

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '112.29'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -5.45%  '
$ws.Range('E31').Style = "Normal"

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.715'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('E32').Style = "Normal"

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.9087'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -7.14%  '
$ws.Range('E33').Style = "Normal"

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.07870'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -5.41%  '
$ws.Range('E34').Style = "Normal"

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '8.358'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -8.96%  '
$ws.Range('E35').Style = "Normal"

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.458'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -8.89%  '
$ws.Range('E36').Style = "Normal"

# Row 37
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '11.04'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.21%  '
$ws.Range('E37').Style = "Normal"

# Row 38
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.008'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('E38').Style = "Normal"

# Row 39
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.848'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -4.71%  '
$ws.Range('E39').Style = "Normal"

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.05666'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -5.26%  '
$ws.Range('E40').Style = "Normal"

# Row 41
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1940'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -4.53%  '
$ws.Range('E41').Style = "Normal"

# Row 42
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.02054'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -8.92%  '
$ws.Range('E42').Style = "Normal"

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.129'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -6.79%  '
$ws.Range('E43').Style = "Normal"

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5459'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -5.95%  '
$ws.Range('E44').Style = "Normal"

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '12.59'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -2.59%  '
$ws.Range('E45').Style = "Normal"

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.620'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -3.17%  '
$ws.Range('E46').Style = "Normal"

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5268'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -5.17%  '
$ws.Range('E47').Style = "Normal"

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.806'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -4.55%  '
$ws.Range('E48').Style = "Normal"

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.105'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('E49').Style = "Normal"

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '111.02'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -5.77%  '
$ws.Range('E50').Style = "Normal"

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06440'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -5.53%  '
$ws.Range('E51').Style = "Normal"

